# Insert a new weekly price record as row 584 in the Espinaca (spinach)
# price sheet. All existing rows from 584 onward shift down by one
# (584->585, 585->586, ..., 685->686), and the new row is populated with
# the latest week's data (date serial 44951 = 2023-01-25).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 584..685 down to 585..686, leaving row 584 empty (Excel
# copies the formatting of the row above, same as a manual "Insert" of a
# whole row).
$ws.Rows.Item(584).Insert()

# Populate the newly inserted row 584 with the new record.
$ws.Cells.Item(584, 1).Value  = 6
$ws.Cells.Item(584, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(584, 3).Value  = "Metropolitana"
$ws.Cells.Item(584, 4).Value  = 44951
$ws.Cells.Item(584, 5).Value  = 13
$ws.Cells.Item(584, 6).Value  = 100112012
$ws.Cells.Item(584, 7).Value  = "Espinaca"
$ws.Cells.Item(584, 8).Value  = "Sin especificar"
$ws.Cells.Item(584, 9).Value  = "Primera"
$ws.Cells.Item(584, 10).Value = 580
$ws.Cells.Item(584, 11).Value = 6000
$ws.Cells.Item(584, 12).Value = 6500
$ws.Cells.Item(584, 13).Value = 6216
$ws.Cells.Item(584, 14).Value = "`$/cuna 10 kilos"
$ws.Cells.Item(584, 15).Value = "Región Metropolitana"
$ws.Cells.Item(584, 16).Value = 622
$ws.Cells.Item(584, 17).Value = 10
$ws.Cells.Item(584, 18).Value = "Hortaliza"
